$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(7).Insert()
$ws.Range("A6:P6").Copy()
$ws.Range("A7:P7").PasteSpecial(-4122)
Write-Host "done"
